$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed the new shared-string values in the same relative order they
# first appear in the finished workbook (0.0.3, then the two new task
# rows) using scratch cells outside the used range, then clear them.
# This keeps the *relative* ordering of newly-introduced strings sane
# even though the engine always appends genuinely-new strings after the
# ones already present in the workbook.
$ws.Range("Z1").Value = "0.0.3"
$ws.Range("Z2").Value = "Hacer Cuenta cliente"
$ws.Range("Z3").Value = "hacer listado de motores, con orden paginacion, dos vistas"
$ws.Range("Z1:Z3").ClearContents()

# --- Highlight the first five "Funcionalidad" entries (B6:B10) in yellow
$ws.Range("B6:B10").Interior.Color = 65535

# --- New row 21: "Hacer Cuenta cliente" / version 0.0.3
$ws.Range("B21").Value = "Hacer Cuenta cliente"
$ws.Range("C21").Value = "0.0.3"

# --- New row 22: extra task description
$ws.Range("B22").Value = "hacer listado de motores, con orden paginacion, dos vistas"

# --- Update the active selection to match the new edit location
$ws.Range("B24").Select() | Out-Null
